$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base pattern of (regcntr_id, machine_id) pairs that repeats every 9 rows
$pairs = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$startRow = 102
$deviceIdStart = 3000121
$rowCount = 45

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $pair = $pairs[$i % 9]
    $regcntrId = $pair[0]
    $machineId = $pair[1]
    $deviceId = $deviceIdStart + $i

    $ws.Cells.Item($row, 1).Value = $regcntrId
    $ws.Cells.Item($row, 2).Value = $machineId
    $ws.Cells.Item($row, 3).Value = $deviceId
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Scroll the view down and select the newly added block, matching the
# author's final on-screen state when they saved the file.
$ws.Range("A102:G146").Select()
$excel.ActiveWindow.ScrollRow = 129

# Page setup was touched (orientation explicitly set to portrait) when the
# author saved from the Page Setup dialog.
$ws.PageSetup.Orientation = 1
